$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Append the new data point (row 8) below the existing A1:A7 series.
$ws.Range("A8").Value = 20

# Excel advances the active selection to the next empty cell after data entry.
$ws.Range("A9").Select()
